$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text values need to be forced to Text format
# so Excel does not coerce them into actual numbers (losing fixed-decimal formatting),
# then the style is reset back to Normal so no stray style index is left on the cell.
$textCells = @("D4","D5","D6","D9","D10","D11","D14","D20","D21","D22","D23","D25","D27","D28","D30","D31","D32","D36","D37","D38","D39","D40","D41","D43","D44","D45","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row
$ws.Range("D2").Value = '66.085.42'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '3.173.63'
$ws.Range("E3").Value = '  -1.48%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '610.64'
$ws.Range("E5").Value = '  +1.04%  '

$ws.Range("D6").Value = '154.56'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.173.61'
$ws.Range("E8").Value = '  -1.46%  '

$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -1.49%  '

$ws.Range("D11").Value = '5.68'
$ws.Range("E11").Value = '  -7.62%  '

$ws.Range("E12").Value = '  +1.34%  '

$ws.Range("E13").Value = '  -1.66%  '

$ws.Range("D14").Value = '38.42'
$ws.Range("E14").Value = '  -2.86%  '

$ws.Range("D15").Value = '3.691.04'
$ws.Range("E15").Value = '  -1.50%  '

$ws.Range("D16").Value = '66.124.84'
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("E17").Value = '  -1.68%  '

$ws.Range("D18").Value = '3.169.86'
$ws.Range("E18").Value = '  -1.90%  '

$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("D20").Value = '511.24'
$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("D21").Value = '15.41'
$ws.Range("E21").Value = '  -0.53%  '

$ws.Range("D22").Value = '0.731'
$ws.Range("E22").Value = '  -1.36%  '

$ws.Range("D23").Value = '8.04'
$ws.Range("E23").Value = '  -0.99%  '

$ws.Range("E24").Value = '  -3.39%  '

$ws.Range("D25").Value = '84.65'
$ws.Range("E25").Value = '  -0.74%  '

$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("D27").Value = '3.02'
$ws.Range("E27").Value = '  -0.25%  '

$ws.Range("D28").Value = '9.14'
$ws.Range("E28").Value = '  -1.30%  '

$ws.Range("E29").Value = '  +3.76%  '

$ws.Range("D30").Value = '3.00'
$ws.Range("E30").Value = '  +4.07%  '

$ws.Range("D31").Value = '7.16'
$ws.Range("E31").Value = '  +4.14%  '

$ws.Range("D32").Value = '28.01'
$ws.Range("E32").Value = '  -0.91%  '

$ws.Range("E33").Value = '  +0.27%  '

$ws.Range("E34").Value = '  -1.96%  '

$ws.Range("E35").Value = '  -1.65%  '

$ws.Range("D36").Value = '502.62'
$ws.Range("E36").Value = '  +3.76%  '

$ws.Range("D37").Value = '55.03'
$ws.Range("E37").Value = '  -0.49%  '

$ws.Range("D38").Value = '0.0884'
$ws.Range("E38").Value = '  -2.59%  '

$ws.Range("D39").Value = '0.0421'
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("D40").Value = '0.128'
$ws.Range("E40").Value = '  +6.62%  '

$ws.Range("D41").Value = '8.79'
$ws.Range("E41").Value = '  -1.91%  '

$ws.Range("D42").Value = '0.0₃0682'
$ws.Range("E42").Value = '  +5.06%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.84'
$ws.Range("E43").Value = '  -4.63%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '0.298'
$ws.Range("E44").Value = '  -1.61%  '

$ws.Range("D45").Value = '2.44'
$ws.Range("E45").Value = '  -0.59%  '

$ws.Range("D46").Value = '2.824.04'

$ws.Range("D47").Value = '28.12'
$ws.Range("E47").Value = '  -2.39%  '

$ws.Range("D48").Value = '2.37'
$ws.Range("E48").Value = '  +2.05%  '

$ws.Range("E49").Value = '  -0.10%  '

$ws.Range("E50").Value = '  +0.28%  '

$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").Value = '2.61'
$ws.Range("E51").Value = '  +6.95%  '

# Reset style back to Normal on the text-protected cells (keeps the text value,
# drops the temporary Text number-format so cell styling matches the original).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Applied crypto price/volume updates"